$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 16 entry: "Reporte de pagos" with amount 900
$ws.Range("B16").Value = "Reporte de pagos"
$ws.Range("C16").Value = 900

# Update sheet view: clear topLeftCell scroll position and move selection to C17
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C17").Select()
